$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto price / 1h-volume table with the latest scrape values,
# and restore WrappedEther/Polkadot (rows 18-19) and Hedera/FirstDigitalUSD
# (rows 33-34) to their new ranking order.
#
# Columns D (Price) and E (Volume(1h)) hold plain text in this sheet, not
# real numbers (prices use "." as a thousands separator, e.g. "61.778.06",
# and volumes are zero-padded percent strings like "  -1.30%  "). Most new
# values are unambiguous text already (multiple dots, %, letters), but a
# few D-column prices look like plain decimals (e.g. "6.52", "0.0000110")
# and Excel would silently reinterpret them as numbers -- losing the
# original text formatting (e.g. trailing zeros). Force those specific
# cells to text format first so they round-trip exactly.

$ws.Range('D2').Value = '61.778.06'
$ws.Range('E2').Value = '  -1.30%  '
$ws.Range('D3').Value = '2.888.89'
$ws.Range('E3').Value = '  -2.13%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '566.58'
$ws.Range('E5').Value = '  -3.86%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '142.94'
$ws.Range('E6').Value = '  -2.64%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  -0.95%  '
$ws.Range('D9').Value = '2.885.31'
$ws.Range('E9').Value = '  -2.17%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '6.96'
$ws.Range('E10').Value = '  +0.11%  '
$ws.Range('E11').Value = '  -1.47%  '
$ws.Range('E12').Value = '  -1.27%  '
$ws.Range('E13').Value = '  -0.39%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '31.88'
$ws.Range('E14').Value = '  -0.84%  '
$ws.Range('E15').Value = '  -0.05%  '
$ws.Range('D16').Value = '3.368.27'
$ws.Range('E16').Value = '  -2.14%  '
$ws.Range('D17').Value = '61.721.91'
$ws.Range('E17').Value = '  -1.42%  '
$ws.Range('B18').Value = 'Polkadot'
$ws.Range('C18').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '6.52'
$ws.Range('E18').Value = '  -1.88%  '
$ws.Range('B19').Value = 'WrappedEther'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D19').Value = '2.891.25'
$ws.Range('E19').Value = '  -2.05%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '428.78'
$ws.Range('E20').Value = '  -1.26%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '12.98'
$ws.Range('E21').Value = '  -3.67%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.653'
$ws.Range('E22').Value = '  -1.19%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '6.86'
$ws.Range('E23').Value = '  -1.34%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '78.90'
$ws.Range('E24').Value = '  -1.49%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '12.03'
$ws.Range('E25').Value = '  +1.13%  '
$ws.Range('E26').Value = '  -9.86%  '
$ws.Range('E27').Value = '  +0.04%  '
$ws.Range('E28').Value = '  -3.35%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.0000110'
$ws.Range('E29').Value = '  +9.47%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '6.94'
$ws.Range('E30').Value = '  -3.21%  '
$ws.Range('E31').Value = '  -3.29%  '
$ws.Range('E32').Value = '  -7.06%  '
$ws.Range('B33').Value = 'FirstDigitalUSD'
$ws.Range('C33').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.00'
$ws.Range('E33').Value = '  -0.10%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.107'
$ws.Range('E34').Value = '  -1.74%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '25.58'
$ws.Range('E35').Value = '  -2.12%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.948'
$ws.Range('E36').Value = '  -4.16%  '
$ws.Range('E37').Value = '  -3.27%  '
$ws.Range('E38').Value = '  -1.74%  '
$ws.Range('E39').Value = '  -6.55%  '
$ws.Range('E40').Value = '  -4.72%  '
$ws.Range('E41').Value = '  +0.56%  '
$ws.Range('E42').Value = '  -2.34%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '39.62'
$ws.Range('E43').Value = '  +1.38%  '
$ws.Range('E44').Value = '  -2.21%  '
$ws.Range('D45').Value = '2.688.08'
$ws.Range('E45').Value = '  +0.22%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0336'
$ws.Range('E46').Value = '  +0.20%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '131.31'
$ws.Range('E47').Value = '  -2.50%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '344.41'
$ws.Range('E48').Value = '  -2.31%  '
$ws.Range('E49').Value = '  -0.04%  '
$ws.Range('E50').Value = '  -1.32%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '21.54'
$ws.Range('E51').Value = '  -4.19%  '

Write-Output "Applied 87 cell updates"
